$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update from_bus / to_bus / in_service values for existing rows 8-15
# (the "name" labels in column B are unchanged).
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

$ws.Range("C9").Value = 16

$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

$ws.Range("C12").Value = 10
$ws.Range("E12").Value = $true

$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11

# New row 16: line7
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "line7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# New row 17: line8
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "line8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false

# Column A on the new rows uses the same (bold/bordered/centered) style as
# the rest of column A in the data rows.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
